$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The exercise's "h" (step) constant in $D$7 is now subtracted instead of
# added when building the finite-difference derivative f'(x); and the
# sample point C5 moves from 1 to 0.
$ws.Range("C5").Value = 0

$ws.Range("E5").Formula = '=(($C$2-C5-$D$7)^2+($D$2-C5-$D$7)^2+($E$2-C5-$D$7)^2-D5)/$D$7'
$ws.Range("D10").Formula = '=(($C$2-B10-$D$7)^2+($D$2-B10-$D$7)^2+($E$2-B10-$D$7)^2-C10)/$D$7'
$ws.Range("D11:D30").Formula = '=(($C$2-B11-$D$7)^2+($D$2-B11-$D$7)^2+($E$2-B11-$D$7)^2-C11)/$D$7'

# Selection moved one row down.
$ws.Range("G6").Select()
